# Roll the 90-day GSC export window forward by two days:
#   - drop 2025-09-01 and 2025-09-02 (and their data row)
#   - append 2025-11-30 and 2025-12-01 (Non-HTTPS = 0, HTTPS = 0, matching
#     the trailing zeros already present at the end of the series)
# Every other row's "HTTPS URLs" count shifts up by two rows as a side
# effect of the deletion, which is exactly what the source diff shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Drop the two oldest dates (rows 2 and 3 -> 2025-09-01 / 2025-09-02);
# this shifts every row below up by two, which is what re-bases the
# "HTTPS URLs" column onto the new date alignment.
$ws.Rows("2:3").Delete()

# Append the two new trailing dates with their data. Cells are formatted
# as text first so Excel stores the literal "yyyy-MM-dd" string instead of
# auto-converting it to a date serial, then the format is cleared so the
# cell's style matches the plain (unstyled) cells around it.
$ws.Range("A90").NumberFormat = "@"
$ws.Range("A90").Value = "2025-11-30"
$ws.Range("A90").ClearFormats()
$ws.Range("B90").Value = 0
$ws.Range("C90").Value = 0

$ws.Range("A91").NumberFormat = "@"
$ws.Range("A91").Value = "2025-12-01"
$ws.Range("A91").ClearFormats()
$ws.Range("B91").Value = 0
$ws.Range("C91").Value = 0
